$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Sum / Average / Median / Max / Min figures (row r corresponds to
# original source rows 2-9 of the sheet).
$data = @{
    2 = @(-58.93939971923828, -0.0442, -0.07000000029802322, 0.5333999991416931, -0.1829999983310699)
    3 = @(-162.0104064941406, -0.1296, -0.183, 0.6333000063896179, -0.1829999983310699)
    4 = @(-130.0059051513672, -0.1519, -0.183, 0.5620999932289124, -0.1829999983310699)
    5 = @(208.3249053955078, 0.2198, 0.2823, 0.8169999718666077, -0.1829999983310699)
    6 = @(216.2174072265625, 0.1711, 0.1224, 0.7531999945640564, -0.1829999983310699)
    7 = @(52.70330047607422, 0.0531, 0.02869999967515469, 0.5097000002861023, -0.1485999971628189)
    8 = @(69.60199737548828, 0.0624, 0.0412, 0.4009999930858612, -0.1485999971628189)
    9 = @(195.8919219970703, 0.0249, 0, 0.8169999718666077, -0.1829999983310699)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
